$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Segment names and values (PercActivations, PercSegmentAreas)
$data = @(
    @("background", 0.2219194031031927, 0.3365718514973398),
    @("back_bumper", 0.006806044373858539, 0.006693558601720179),
    @("back_glass", 0.01495318657512381, 0.01923960726569185),
    @("back_left_door", 0.004766413912031224, 0.004202819515640954),
    @("back_left_light", 0.005136159929754412, 0.004407771227443971),
    @("back_right_door", 0.006469493272063987, 0.00561351725119662),
    @("back_right_light", 0.006775173181008223, 0.005856581263861952),
    @("front_bumper", 0.3017020994991718, 0.2289109289474905),
    @("front_glass", 0.1491447472621007, 0.1936832406181643),
    @("front_left_door", 0.002538274891617849, 0.001816508462775329),
    @("front_left_light", 0.02299154253470514, 0.0146677320671856),
    @("front_right_door", 0.002427235949421697, 0.00183439196152098),
    @("front_right_light", 0.0181456157772409, 0.01419983133905514),
    @("hood", 0.2251117279708277, 0.1511630354932985),
    @("left_mirror", 0.003371874176780485, 0.003833988049159126),
    @("right_mirror", 0.002150726685720349, 0.003205786059830067),
    @("tailgate", 0.0002404645567990321, 0.0001659201237684006),
    @("trunk", 0.005327969854529005, 0.003889687883140958),
    @("wheel", 0.0000218463168277382, 0.0000432423717161807)
)

# Insert a new column before column B. The old B/C columns (values) shift to
# C/D, keeping their existing (unstyled) formatting. Column A (segment names,
# styled) stays put and will be repurposed to hold the numeric index.
$ws.Columns.Item(2).Insert()

# New column B (segments) should carry no special style, matching the
# un-styled value columns from before the insert.
$ws.Range("B1:B20").ClearFormats()

# Header row
$ws.Cells.Item(1, 2).Value = "segments"
$ws.Cells.Item(1, 3).Copy()
$ws.Cells.Item(1, 2).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows: column A becomes a numeric 0-based index (keeps the old header
# style), column B gets the segment name.
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i
    $ws.Cells.Item($row, 2).Value = $data[$i][0]
}

$ws.Range("A1:D20").Value2
